$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-9: Genre, latitude (B), longitude (C)
$data = @(
    @("gospel",      33.80775722088888, 31.63713273881481),
    @("letter",       34.53737994731184, 29.36316048422938),
    @("historical",   31.9524206297491,  36.09865291409798),
    @("lyric",        31.05265941960786, 35.46695629150328),
    @("wisdom",       27.77928134285714, 35.48837334285714),
    @("prophecy",     32.04909800329218, 35.39011386927299),
    @("apocalyptic",  32.24773232303031, 35.54896866929293),
    @("law",          31.65596558106063, 35.77621887525253)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
